$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Update the caption of the existing cross-validation block
#    (G8) to rich text: "Tabella REPORT5FOLD: Report della "
#    + "5 fold cross" (red) + " validation del mo-dello SVM con
#    ProfilingUD " (black).
#    (Done first so this shared string keeps its original slot.)
# -----------------------------------------------------------------
$full = "Tabella REPORT5FOLD: Report della 5 fold cross validation del mo-dello SVM con ProfilingUD "
$ws.Range("G8").Value = $full

$run2 = $ws.Range("G8").Characters(35, 12)
$run2.Font.Color = 255
$run2.Font.Name = "Times New Roman"
$run2.Font.Size = 11

$run3 = $ws.Range("G8").Characters(47, 45)
$run3.Font.Color = 0
$run3.Font.Name = "Times New Roman"
$run3.Font.Size = 11

# -----------------------------------------------------------------
# 2) Add the new "final test-set report" block (M2:Q8), mirroring
#    the existing SVM ProfilingUD per-genre block in G2:K8.
# -----------------------------------------------------------------
$ws.Range("G2:K8").Copy($ws.Range("M2"))
$excel.CutCopyMode = 0

# New caption for the pasted block (replaces the copied "Tabella XXXX..." text)
$ws.Range("M8").Value = "Report finale sul test set di SVM ProfilingUD "

# New metrics for the final test-set report (row 4: CH only).
# "0.5729"/"0.57" must land as *text* (like the existing "0.51"/"0.50"
# neighbours) without Excel re-parsing them into numbers and without
# leaving the cell's number format changed, so stage them in a scratch
# cell formatted as Text and paste-values them across.
$ws.Range("Y1").NumberFormat = "@"
$ws.Range("Y1").Value = "0.5729"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "0.57"

$ws.Range("Y1").Copy()
$ws.Range("N4").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("Y2").Copy()
$ws.Range("O4").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("Y1:Y2").Clear()

$ws.Range("P4").Value = 199
$ws.Range("Q4").Value = "[[69, 31], [54, 45]]"

# Give column Q a bit more width (matches the custom width added for the new caption column)
$ws.Columns("Q").ColumnWidth = 14.86328125

# -----------------------------------------------------------------
# 3) Update the view: scroll so column E is the leftmost visible
#    column, and select the newly added block.
# -----------------------------------------------------------------
$ws.Range("M2:Q8").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5

# -----------------------------------------------------------------
# 4) Tidy up leftover stray formatting at the bottom of the sheet
#    that is no longer needed (rows 24:25 were blank formatted
#    rows below the last real content row).
# -----------------------------------------------------------------
$ws.Rows("24:25").Delete()
